# Sample Project / Main.xlsx - "Rules" sheet, rule R40 (row 11):
# B11 held the shared string "R40"; it now holds the literal text "1"
# (kept as text, not converted to a number), while every other
# attribute of the cell (style s="23", row height, etc.) stays the same.
#
# A plain `$cell.Value = "1"` (or `.Value = "'1"`) makes Excel's input
# parser treat the text as a number / apply a quote-prefix style, which
# bumps the cell's style id. Routing the literal text through a
# formula ("=""1""") and pasting only the *values* keeps the cell's
# existing formatting/style untouched while still storing "1" as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Cells.Item(11, 2)   # B11

$helper = $ws.Cells.Item(1, 20)   # scratch cell, well outside the used range
$helper.Formula = '="1"'
$helper.Copy()
$target.PasteSpecial(-4163)       # xlPasteValues - value only, formatting untouched
$helper.Clear()
$excel.CutCopyMode = $false
